$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignInContactUs")

# -------------------------------------------------------------------------
# Stash a copy of the existing "Hyperlink" cell format (style) from B3 onto
# a scratch cell far outside the used range, so we can restore it later
# after re-creating the hyperlinks (Hyperlinks.Add re-applies its own
# built-in Hyperlink style/font, which would otherwise leave the cells
# pointing at a duplicated style record instead of the original one).
# -------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# Rebuild the hyperlinks collection (new IDs/Project/Suite/Section/TestCase
# values + refreshed addresses) in the same order as before so relationship
# ids line up the same way on export.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:abc@xyz.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:abc@xyz.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Test@123", "", "", "Test@123")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:Test@123", "", "", "Test@123")
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:abc@123.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:abc@234.com")

# Restore the original Hyperlink cell style on every cell the rebuild above
# touched (it otherwise swaps in a freshly minted - but visually identical -
# style record).
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# -------------------------------------------------------------------------
# Updated IDs / credentials / email used by the automation suite.
# -------------------------------------------------------------------------
$ws.Range("A2").Value = "testcase-001"
$ws.Range("A3").Value = "testcase-001"
$ws.Range("B2").Value = "alistair.zhu@laserfiche.com"
$ws.Range("B3").Value = "alistair.zhu@laserfiche.com"
$ws.Range("C2").Value = "19Unipas91"
$ws.Range("C3").Value = "19Unipas91"
$ws.Range("F2").Value = "alistair.zhu@laserfiche.com"
$ws.Range("F3").Value = "alistair.zhu@laserfiche.com"

# Column widths widened now that the new values are longer (auto-fit no
# longer leaves these at their old "best fit" widths).
$ws.Columns.Item(2).ColumnWidth = 34.27396
$ws.Columns.Item(3).ColumnWidth = 18.9888
$ws.Columns.Item(4).ColumnWidth = 17.70365
$ws.Columns.Item(6).ColumnWidth = 32.27396
$ws.Columns.Item(7).ColumnWidth = 21.58345
$ws.Columns.Item(8).ColumnWidth = 24.13333

# Selection moved to A3.
$ws.Range("A3").Select()
